$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A5").Value = "13.4.21"
$ws.Range("C5").Value = 24
$ws.Range("F5").Value = "left"
$ws.Range("G5").Value = "N"
$ws.Range("I5").Value = "F"
$ws.Range("J5").Value = "Y"
$ws.Range("L5").Value = "Y"
$ws.Range("M5").Value = "Y"

$ws.Hyperlinks.Add($ws.Range("N5"), "mailto:netayellin97@walla.com", "", "", "netayellin97@walla.com")

$ws.Range("T5").Select()
$excel.ActiveWindow.ScrollColumn = 17
$excel.ActiveWindow.ScrollRow = 1
